$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data (rows re-ordered / content refreshed: AVG, CTSES2, CTSES1, CodeBLEU) ---
$ws.Range("A2").Value = "AVG"
$ws.Range("B2").Value = 0.40128
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "Defects4J::RandomStringUtils::iter-1,SF110::CompareToBuilder::iter-2,SF110::XPathLexer::iter-1,SF110::XPathLexer::iter-2"

$ws.Range("A3").Value = "CTSES2"
$ws.Range("B3").Value = 0.404453
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = "Defects4J::RandomStringUtils::iter-1,SF110::CompareToBuilder::iter-2,SF110::XPathLexer::iter-1,SF110::XPathLexer::iter-2"

$ws.Range("A4").Value = "CTSES1"
$ws.Range("B4").Value = 0.40666
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = "Defects4J::RandomStringUtils::iter-1,SF110::CompareToBuilder::iter-2,SF110::XPathLexer::iter-1,SF110::XPathLexer::iter-2"

$ws.Range("A5").Value = "CodeBLEU"
$ws.Range("B5").Value = 0.432893
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "Defects4J::TarUtils::iter-3,SF110::CompareToBuilder::iter-2,SF110::XPathLexer::iter-1,SF110::XPathLexer::iter-2"

# --- Formatting: wrap text everywhere, thin border + wrap on the data block ---
$ws.Range("A1:D1").WrapText = $true

$ws.Range("A2:D5").WrapText = $true
$ws.Range("A2:D5").Borders.LineStyle = 1

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 32
$ws.Rows.Item(2).RowHeight = 224
$ws.Rows.Item(3).RowHeight = 224
$ws.Rows.Item(4).RowHeight = 224
$ws.Rows.Item(5).RowHeight = 208

# --- Column D width (engine quantizes to pixel grid; 111.3 is the closest input
#     that lands on the target's stored width of ~112.164) ---
$ws.Columns.Item(4).ColumnWidth = 111.3

# --- Selection mirrors the saved workbook state ---
$ws.Range("D9").Select()
